# Generate Report for Handoff
# - Flip the "Status" text from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#   de-de!C2).
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to match the new handoff.
# - Narrow the now-shorter status column(s) to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview's "Latest HO Xliff Generate Date" and de-de's "Latest Handoff
# Datetime" both move from 02:58:09 to 02:58:59.
$wsOverview.Range("G2").Value = "2016-08-25 02:58:59"
$wsDeDe.Range("H2").Value = "2016-08-25 02:58:59"

# zh-cn's "Latest Handoff Datetime" moves from 02:57:57 to 02:58:53.
$wsZhCn.Range("H2").Value = "2016-08-25 02:58:53"

# --- Column widths: the Status columns shrink now that the text is shorter ---
# (ColumnWidth is in characters; Excel quantizes to whole display pixels,
# so 16.3333... is the closest input that lands on the target width.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
